# Append the new daily row (row 91) to the tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a date-looking string ("2025/10/11") that must be stored as
# literal text (matching the rest of the column), not auto-converted to a
# serial date number. Temporarily mark the cell as Text, write the value,
# then clear the format again so the cell is left with the sheet's default
# (unstyled) formatting, same as every other data row.
$ws.Range("A91").NumberFormat = "@"
$ws.Range("A91").Value = "2025/10/11"
$ws.Range("A91").ClearFormats()

$ws.Range("B91").Value = "土"
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 37
